$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 50.855869642042443
$ws.Range("B2").Value = 0.97924340846706515
$ws.Range("C2").Value = 0.16478512100587317
$ws.Range("D2").Value = -0.11801784131598221

$ws.Range("A3").Value = -8.3524547089766639
$ws.Range("B3").Value = 0.18860496310208769
$ws.Range("C3").Value = -0.5275597259011382
$ws.Range("D3").Value = 0.82831691006545061

$ws.Range("A4").Value = -74.759483710354601
$ws.Range("B4").Value = 0.074232842240242666
$ws.Range("C4").Value = -0.83338262491018833
$ws.Range("D4").Value = -0.54768867582846692

$ws.Range("G4").Select()
